$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D46").Value = "로봇 vs 복강경 위절제술"
$ws.Range("E46").Value = "https://bioinformaticsandme.tistory.com/505"

$ws.Range("D50").Value = "기계학습으로 합금 발견"
$ws.Range("E50").Value = "http://incredible.egloos.com/7573300"

$ws.Range("D51").Value = "[Excel 365] 엑셀에서 범위를 선택할 때 $를 사용하는 경우"
$ws.Range("E51").Value = "https://bskyvision.com/entry/Excel-365-%EC%97%91%EC%85%80%EC%97%90%EC%84%9C-%EB%B2%94%EC%9C%84%EB%A5%BC-%EC%84%A0%ED%83%9D%ED%95%A0-%EB%95%8C-%EB%A5%BC-%EC%82%AC%EC%9A%A9%ED%95%98%EB%8A%94-%EA%B2%BD%EC%9A%B0"
